$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.228.02'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.71%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.914.19'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.19%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9996'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.23%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.8230'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +5.04%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '244.42'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.09%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.9992'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.20%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3262'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +3.40%  '

$ws.Range('E9').Value = '  +4.46%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07120'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.25%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.08085'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.62%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.7801'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +5.14%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.919.18'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.49%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.361'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.93%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '94.55'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.72%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '30.250.20'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.75%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.35'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.77%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.009'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +2.27%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '249.12'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.27%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000007845'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.23%  '

$ws.Range('E21').Value = '  -0.03%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9996'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.18%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.630'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +11.03%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.1683'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +22.60%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.460'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.91%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '168.35'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.42%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.13'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.04%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.122'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +4.58%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.368'
$ws.Range('D29').Style = "Normal"

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.533'
$ws.Range('D30').Style = "Normal"

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.345'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.64%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.05681'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +3.80%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.135'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.30%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.291'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.87%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7403'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.86%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.001'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.04%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.714'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.61%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01940'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.32%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.812'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +0.66%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.4485'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +1.63%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '73.99'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.73%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.985'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.38%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.935'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +3.12%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.8495'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.44%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.9995'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.19%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.036.59'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +5.76%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '102.97'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.61%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.958'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +2.15%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.622'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.40%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.088.14'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.46%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.562'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +4.86%  '
